$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vessel_sf")

# Remove the "Turntable outer diameter [m]" (row 9) and "Dredge depth [m]" (row 10)
# rows from the vessel safety-factor table. The remaining rows shift up but keep
# their original id numbers (6, 9, 10, 11, 12, ...).
$ws.Rows(9).Delete()
$ws.Rows(9).Delete()

# Make vessel_sf the active sheet/tab, with rows 9:10 (now AH winch / AH drum rows)
# selected as whole rows.
$ws.Activate()
$ws.Range("A9:A10").EntireRow.Select()
